$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.469.49"
$ws.Range("E2").Value = "  +2.64%  "
$ws.Range("D3").Value = "'1.721.95"
$ws.Range("E3").Value = "  +2.45%  "
$ws.Range("D4").Value = "'1.006"
$ws.Range("E4").Value = "  +0.98%  "
$ws.Range("D5").Value = "'225.01"
$ws.Range("E5").Value = "  -1.10%  "
$ws.Range("D6").Value = "'0.5337"
$ws.Range("E6").Value = "  +0.57%  "
$ws.Range("E7").Value = "  +0.88%  "
$ws.Range("D8").Value = "'0.2655"
$ws.Range("E8").Value = "  -0.58%  "
$ws.Range("D9").Value = "'0.06587"
$ws.Range("E9").Value = "  +3.03%  "
$ws.Range("D10").Value = "'21.40"
$ws.Range("E10").Value = "  +2.18%  "
$ws.Range("D11").Value = "'0.07671"
$ws.Range("E11").Value = "  -0.33%  "
$ws.Range("D12").Value = "'4.594"
$ws.Range("E12").Value = "  -1.78%  "
$ws.Range("D13").Value = "'1.727.36"
$ws.Range("E13").Value = "  +2.93%  "
$ws.Range("D14").Value = "'1.960.95"
$ws.Range("E14").Value = "  +3.36%  "
$ws.Range("D15").Value = "'0.5785"
$ws.Range("E15").Value = "  +0.72%  "
$ws.Range("D16").Value = "'0.0₅8273"
$ws.Range("E16").Value = "  +0.11%  "
$ws.Range("D17").Value = "'67.71"
$ws.Range("E17").Value = "  +1.92%  "
$ws.Range("D18").Value = "'27.488.37"
$ws.Range("E18").Value = "  +3.85%  "
$ws.Range("D19").Value = "'217.55"
$ws.Range("E19").Value = "  +9.05%  "
$ws.Range("D20").Value = "'1.005"
$ws.Range("E20").Value = "  +0.42%  "
$ws.Range("D21").Value = "'4.731"
$ws.Range("E21").Value = "  +0.61%  "
$ws.Range("D22").Value = "'10.55"
$ws.Range("E22").Value = "  -1.45%  "
$ws.Range("D23").Value = "'5.991"
$ws.Range("E23").Value = "  -0.97%  "
$ws.Range("D24").Value = "'1.006"
$ws.Range("E24").Value = "  +0.92%  "
$ws.Range("D25").Value = "'143.37"
$ws.Range("E25").Value = "  -1.57%  "
$ws.Range("D26").Value = "'1.731"
$ws.Range("E26").Value = "  +9.86%  "
$ws.Range("D27").Value = "'0.1229"
$ws.Range("E27").Value = "  +0.85%  "
$ws.Range("D28").Value = "'7.314"
$ws.Range("E28").Value = "  -0.01%  "
$ws.Range("D29").Value = "'16.45"
$ws.Range("E29").Value = "  +1.27%  "
$ws.Range("D30").Value = "'0.05431"
$ws.Range("E30").Value = "  -3.35%  "
$ws.Range("D31").Value = "'1.298"
$ws.Range("E31").Value = "  -0.42%  "
$ws.Range("D32").Value = "'3.538"
$ws.Range("E32").Value = "  -0.45%  "
$ws.Range("D33").Value = "'3.423"
$ws.Range("E33").Value = "  -0.61%  "
$ws.Range("D34").Value = "'1.645"
$ws.Range("E34").Value = "  +3.53%  "
$ws.Range("D35").Value = "'2.880"
$ws.Range("E35").Value = "  +2.79%  "
$ws.Range("D36").Value = "'0.9553"
$ws.Range("E36").Value = "  -0.38%  "
$ws.Range("D37").Value = "'2.431"
$ws.Range("E37").Value = "  +0.79%  "
$ws.Range("D38").Value = "'0.5905"
$ws.Range("E38").Value = "  +2.40%  "
$ws.Range("E39").Value = "  +2.09%  "
$ws.Range("D40").Value = "'5.901"
$ws.Range("E40").Value = "  +1.23%  "
$ws.Range("D41").Value = "'1.047.15"
$ws.Range("E41").Value = "  +0.57%  "
$ws.Range("D42").Value = "'0.8452"
$ws.Range("E42").Value = "  +1.05%  "
$ws.Range("E43").Value = "  +0.56%  "
$ws.Range("D44").Value = "'101.14"
$ws.Range("E44").Value = "  -0.84%  "
$ws.Range("D45").Value = "'1.867.77"
$ws.Range("E45").Value = "  +3.15%  "
$ws.Range("D46").Value = "'0.0₈115"
$ws.Range("E46").Value = "  +8.00%  "
$ws.Range("D47").Value = "'58.47"
$ws.Range("E47").Value = "  -0.37%  "
$ws.Range("D48").Value = "'0.4511"
$ws.Range("E48").Value = "  +4.30%  "

# Row 49: EnergySwap -> Frax
$ws.Range("B49").Value = "Frax"
$ws.Range("C49").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D49").Value = "'1.004"
$ws.Range("E49").Value = "  -1.09%  "

# Row 50: Frax -> EnergySwap
$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D50").Value = "'8.128"
$ws.Range("E50").Value = "  +0.53%  "

# Row 51
$ws.Range("D51").Value = "'0.06553"
$ws.Range("E51").Value = "  +12.78%  "
